$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace existing values and add new ones in column A
$values = @("cà phê", "cf", "trà", "cà phê sách", "quán ăn", "nhà hàng", "sushi")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Move active cell selection to the row after the last entry, like after typing+Enter
$ws.Range("A8").Select()
